# This script applies the "Deploying to gh-pages" content refresh to the
# StructureDefinition-source-record-id workbook:
#   - bump Version 5.0.0 -> 6.0.0
#   - bump Date to the new publish timestamp
#   - replace the placeholder "Contact" rows with real Publisher/Jurisdiction
#     metadata (Alvearie Team / United States of America)
#   - update the root Extension row's Short/Definition text on the
#     "Elements" sheet to describe the SourceRecordId extension specifically

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# --- Metadata sheet -------------------------------------------------------
# Row 11 is the second, redundant "Contact" / "No display for ContactDetail"
# row; delete it so the table collapses back to one row per property, then
# fill in the real values for Publisher/Jurisdiction.
$meta.Rows.Item(11).Delete()

$meta.Range("B3").Value = "6.0.0"
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"
$meta.Range("B9").Value = "Alvearie Team"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# --- Elements sheet --------------------------------------------------------
# The root Extension row's Short/Definition columns (K/L) previously held the
# generic placeholder text; update them to describe this specific extension.
$elements.Range("K2").Value = "Source Record ID"
$elements.Range("L2").Value = "The ID for a record that the data producer or data integrator extracted knowledge from to produce the data within the FHIR resource or element"
